$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(17).Delete()
